$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '42.735.94'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -1.17%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.370.02'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.83%  '
$ws.Range('E4').Value = '  -0.17%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '332.48'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +5.73%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '101.06'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -7.93%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.638'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.62%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.629'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.60%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '39.99'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -6.90%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0920'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.92%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '8.49'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -4.44%  '
$ws.Range('E13').Value = '  -3.79%  '
$ws.Range('E14').Value = '  +0.10%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '16.47'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.95%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.725.13'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.89%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.370.84'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.12%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '8.04'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +10.54%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '42.667.20'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.26%  '
$ws.Range('E20').Value = '  -1.71%  '
$ws.Range('E21').Value = '  +9.65%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '76.45'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.45%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '270.02'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +6.24%  '
$ws.Range('E24').Value = '  -11.05%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '10.15'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +11.41%  '
$ws.Range('E26').Value = '  -0.09%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.53'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -4.49%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '23.23'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +3.70%  '
$ws.Range('E29').Value = '  -2.67%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '176.60'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.00%  '
$ws.Range('E31').Value = '  -2.52%  '
$ws.Range('E32').Value = '  -2.59%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '35.35'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -10.01%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.15'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +2.05%  '
$ws.Range('E35').Value = '  +0.36%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.62'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -7.02%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.99'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +10.91%  '
$ws.Range('E38').Value = '  -4.71%  '
$ws.Range('B39').Value = 'Kaspa'
$ws.Range('C39').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.105'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.47%  '
$ws.Range('B40').Value = 'NEARProtocol'
$ws.Range('C40').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.82'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -8.05%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.53'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +4.02%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.235'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.69%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '70.21'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -3.78%  '
$ws.Range('E44').Value = '  -0.04%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '118.07'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +6.21%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '91.56'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +31.05%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '11.90'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -7.31%  '
$ws.Range('B48').Value = 'FraxShare'
$ws.Range('C48').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.26'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.61%  '
$ws.Range('B49').Value = 'THORChain'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '5.50'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.38%  '
$ws.Range('E50').Value = '  -2.75%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.571.22'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +5.11%  '
